$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update BOM description label (formerly "Part 50pcs")
$ws.Range("A1").Value = "Part (50pcs)"

# Stencils count updated from 4 to 3
$ws.Range("C16").Value = 3

# Selection state as recorded in the saved workbook
$ws.Activate()
$ws.Range("C17").Select()
